$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column B ("l" = linear flag) to 1 for rows 2-6 (trials now nonlinear)
$ws.Range("B2").Value = 1
$ws.Range("B3").Value = 1
$ws.Range("B4").Value = 1
$ws.Range("B5").Value = 1
$ws.Range("B6").Value = 1

# Update column D ("q" = proximity) to 2 for rows 2-7
$ws.Range("D2").Value = 2
$ws.Range("D3").Value = 2
$ws.Range("D4").Value = 2
$ws.Range("D5").Value = 2
$ws.Range("D6").Value = 2
$ws.Range("D7").Value = 2

# Update the active selection to match the recorded view state
$ws.Range("E17").Select()
